# Edits to boxplots and statistics for authors
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header rename: "transport" -> "group"
$ws.Cells.Item(1, 3).Value = "group"

# Rows where column B == "alltextures" -> "All Textures"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "All Textures"
}

# Rows where column B == "mechanical" -> "Mechanical"
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).Value = "Mechanical"
}
